$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-20 Friday", 2) | Out-Null
$d.Content.Find.Execute("58×90=", $true, $false, $false, $false, $false, $true, 1, $false, "94×29=", 2) | Out-Null
$d.Content.Find.Execute("89×83=", $true, $false, $false, $false, $false, $true, 1, $false, "54×16=", 2) | Out-Null
$d.Content.Find.Execute("72×24=", $true, $false, $false, $false, $false, $true, 1, $false, "29×40=", 2) | Out-Null
$d.Content.Find.Execute("46×15=", $true, $false, $false, $false, $false, $true, 1, $false, "92×41=", 2) | Out-Null
$d.Content.Find.Execute("85×49=", $true, $false, $false, $false, $false, $true, 1, $false, "82×83=", 2) | Out-Null
$d.Content.Find.Execute("54×92=", $true, $false, $false, $false, $false, $true, 1, $false, "85×76=", 2) | Out-Null
$d.Content.Find.Execute("43×38=", $true, $false, $false, $false, $false, $true, 1, $false, "31×73=", 2) | Out-Null
$d.Content.Find.Execute("30×76=", $true, $false, $false, $false, $false, $true, 1, $false, "37×40=", 2) | Out-Null
$d.Content.Find.Execute("96×52=", $true, $false, $false, $false, $false, $true, 1, $false, "69×35=", 2) | Out-Null
$d.Content.Find.Execute("26×81=", $true, $false, $false, $false, $false, $true, 1, $false, "28×50=", 2) | Out-Null
$d.Content.Find.Execute("86×72=", $true, $false, $false, $false, $false, $true, 1, $false, "61×59=", 2) | Out-Null
$d.Content.Find.Execute("59×37=", $true, $false, $false, $false, $false, $true, 1, $false, "83×32=", 2) | Out-Null
$d.Content.Find.Execute("76×24=", $true, $false, $false, $false, $false, $true, 1, $false, "17×83=", 2) | Out-Null
$d.Content.Find.Execute("84×82=", $true, $false, $false, $false, $false, $true, 1, $false, "69×65=", 2) | Out-Null
$d.Content.Find.Execute("76×57=", $true, $false, $false, $false, $false, $true, 1, $false, "57×57=", 2) | Out-Null
$d.Content.Find.Execute("60×65=", $true, $false, $false, $false, $false, $true, 1, $false, "71×49=", 2) | Out-Null
$d.Content.Find.Execute("78×83=", $true, $false, $false, $false, $false, $true, 1, $false, "25×39=", 2) | Out-Null
$d.Content.Find.Execute("42×18=", $true, $false, $false, $false, $false, $true, 1, $false, "95×89=", 2) | Out-Null
$d.Content.Find.Execute("30×47=", $true, $false, $false, $false, $false, $true, 1, $false, "66×42=", 2) | Out-Null
$d.Content.Find.Execute("61×91=", $true, $false, $false, $false, $false, $true, 1, $false, "86×76=", 2) | Out-Null
$d.Content.Find.Execute("30×43=", $true, $false, $false, $false, $false, $true, 1, $false, "65×46=", 2) | Out-Null
$d.Content.Find.Execute("90×85=", $true, $false, $false, $false, $false, $true, 1, $false, "92×28=", 2) | Out-Null
$d.Content.Find.Execute("30×13=", $true, $false, $false, $false, $false, $true, 1, $false, "41×38=", 2) | Out-Null
$d.Content.Find.Execute("30×86=", $true, $false, $false, $false, $false, $true, 1, $false, "54×57=", 2) | Out-Null
$d.Content.Find.Execute("72×20=", $true, $false, $false, $false, $false, $true, 1, $false, "68×19=", 2) | Out-Null
